$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 712204551721.0481
$ws.Range("C3").Value = 242724325270.2232
$ws.Range("C4").Value = 38347168320.21263
$ws.Range("C5").Value = 34385820455.79092
$ws.Range("C6").Value = 26538465526.84064
$ws.Range("C7").Value = 13131836977.62074
$ws.Range("C8").Value = 10570140100.82911
$ws.Range("C9").Value = 9249568473.746784
$ws.Range("C10").Value = 8465833466.668921
$ws.Range("C11").Value = 8301810750.139015
$ws.Range("C12").Value = 7961214387.485915
$ws.Range("C13").Value = 7217284884.089936
$ws.Range("C14").Value = 6794857449.435575
$ws.Range("C15").Value = 5954412669.867033
$ws.Range("C16").Value = 5340065056.106133
$ws.Range("C17").Value = 5094964998.879547
$ws.Range("C18").Value = 4571873154.264347
$ws.Range("C19").Value = 3768600131.609521
$ws.Range("C20").Value = 3614769306.166145
$ws.Range("C21").Value = 3580676608.937095
$ws.Range("C22").Value = 3373783328.34027
$ws.Range("C23").Value = 3009527136.437813
$ws.Range("C24").Value = 2805664665.719321
$ws.Range("C25").Value = 2602614942.655396
$ws.Range("C26").Value = 2499993416.47157
